$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "40.063.52"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  -0.13%  "

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "2.213.03"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  -0.96%  "

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "0.999"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "  -0.21%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "294.38"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -0.22%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "87.97"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +1.39%  "

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.510"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  -1.20%  "

$ws.Range("E8").Value = "  -0.06%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.474"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +0.14%  "

$ws.Range("B10").Value = "Avalanche"
$ws.Range("C10").Value = "https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax"
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "30.04"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -4.06%  "

$ws.Range("B11").Value = "Dogecoin"
$ws.Range("C11").Value = "https://coinranking.com/coin/a91GCGd_u96cF+dogecoin-doge"
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0779"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -1.65%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "49.79"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  +5.88%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.112"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  +2.77%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "6.47"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +0.11%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "2.551.08"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  -0.80%  "

$ws.Range("E16").Value = "  -2.75%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "2.199.15"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  -3.27%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "0.727"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  -0.64%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "39.954.67"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -0.16%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "0.0₃0886"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -0.61%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "11.29"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +3.59%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "5.77"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -0.51%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "65.31"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +0.19%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "237.73"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +0.80%  "

$ws.Range("E25").Value = "  -0.12%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "2.46"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -0.24%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "1.80"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -2.64%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "22.49"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -1.81%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "2.16"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -3.07%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "9.18"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -0.65%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "156.65"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +2.73%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "31.32"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -6.13%  "

$ws.Range("E33").Value = "  -0.10%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "4.90"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +0.40%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.0711"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -0.93%  "

$ws.Range("E36").Value = "  -2.17%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "2.83"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +3.66%  "

$ws.Range("E38").Value = "  +0.49%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.0979"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -2.22%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "15.38"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -6.09%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "1.67"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -1.57%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "2.123.60"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +4.02%  "

$ws.Range("E43").Value = "  -2.57%  "

$ws.Range("E44").Value = "  -2.60%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.0268"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -1.31%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "17.46"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +6.11%  "

$ws.Range("E47").Value = "  -4.99%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "2.65"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +2.77%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "2.420.33"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -0.94%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "1.49"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +2.41%  "

$ws.Range("E51").Value = "  +0.77%  "
